# Apply the recorded edits to param_data_naninovel.xlsx:
#  1. Rename the "Speaker" sheet to "Name".
#  2. Make the renamed "Name" sheet the active sheet/tab, with the
#     selection set to J15 (was B18 before the rename).
#  3. The previously-active sheet ("PackedEffect") is no longer the
#     active tab; its own selection (B5) is left untouched.

$wb = $excel.ActiveWorkbook

# --- 1. Rename "Speaker" -> "Name" -------------------------------------
$speaker = $wb.Worksheets.Item("Speaker")
$speaker.Name = "Name"

# --- 2. Activate "Name" and move its selection to J15 ------------------
$nameSheet = $wb.Worksheets.Item("Name")
$nameSheet.Activate() | Out-Null
$nameSheet.Range("J15").Select() | Out-Null
